# Auto-generated edit script applying the cryptos price-list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.016.73'
$ws.Range("E2").Value = '  +2.93%  '
$ws.Range("D3").Value = '2.951.10'
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'594.93"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = "'148.58"
$ws.Range("E6").Value = '  +2.58%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '2.947.19'
$ws.Range("E8").Value = '  +0.89%  '
$ws.Range("D9").Value = "'0.508"
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("D10").Value = "'7.25"
$ws.Range("E10").Value = '  +4.04%  '
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = '  +7.33%  '
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("E13").Value = '  +5.85%  '
$ws.Range("D14").Value = "'32.87"
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '3.441.67'
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").Value = '62.977.18'
$ws.Range("E17").Value = '  +2.82%  '
$ws.Range("D18").Value = "'6.71"
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").Value = '2.947.54'
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").Value = "'443.18"
$ws.Range("E20").Value = '  +2.91%  '
$ws.Range("D21").Value = "'13.45"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = "'11.30"
$ws.Range("E24").Value = '  +4.27%  '
$ws.Range("D25").Value = "'81.14"
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("D27").Value = "'11.83"
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("D30").Value = "'7.23"
$ws.Range("E30").Value = '  +4.94%  '
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("E32").Value = '  +17.17%  '
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("D34").Value = "'26.44"
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.20%  '
$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("D37").Value = "'3.16"
$ws.Range("E37").Value = '  +6.51%  '
$ws.Range("D38").Value = "'5.62"
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = "'2.05"
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = "'49.71"
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("D43").Value = "'0.282"
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("D44").Value = "'39.03"
$ws.Range("E44").Value = '  -7.25%  '
$ws.Range("D45").Value = "'135.38"
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("D46").Value = '2.695.06'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = "'0.0338"
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("D48").Value = "'362.73"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").Value = "'22.79"
$ws.Range("E51").Value = '  -2.93%  '
